$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44355
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 7000
$ws.Range("P2").Value = 194

# Row 3
$ws.Range("D3").Value = 44372
$ws.Range("J3").Value = 150

# Row 4
$ws.Range("D4").Value = 44701
$ws.Range("J4").Value = 200
$ws.Range("O4").Value = 'Región del Maule'

# Row 6
$ws.Range("D6").Value = 44369
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("N6").Value = '$/caja 20 docenas'
$ws.Range("P6").Value = 7000
$ws.Range("Q6").Value = 1

# Row 7
$ws.Range("D7").Value = 44342
$ws.Range("J7").Value = 150

# Row 8
$ws.Range("D8").Value = 44376
$ws.Range("K8").Value = 6500
$ws.Range("L8").Value = 6500
$ws.Range("M8").Value = 6500
$ws.Range("P8").Value = 181

# Row 9
$ws.Range("D9").Value = 44707
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 9000
$ws.Range("O9").Value = 'Región Metropolitana'
$ws.Range("P9").Value = 250

# Row 10
$ws.Range("D10").Value = 44706
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 9000
$ws.Range("O10").Value = 'Región Metropolitana'
$ws.Range("P10").Value = 250

# Row 11
$ws.Range("D11").Value = 44357
$ws.Range("J11").Value = 150
$ws.Range("N11").Value = '$/caja 20 docenas'
$ws.Range("O11").Value = 'Región del Maule'
$ws.Range("P11").Value = 6500
$ws.Range("Q11").Value = 1

# Row 12
$ws.Range("D12").Value = 44386
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 6500
$ws.Range("M12").Value = 6500
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 181

# Row 14
$ws.Range("D14").Value = 44340
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("N14").Value = '$/caja 36 atados'
$ws.Range("P14").Value = 194
$ws.Range("Q14").Value = 36

# Row 15
$ws.Range("D15").Value = 44348
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("O15").Value = 'Región del Maule'
$ws.Range("P15").Value = 194

# Row 16
$ws.Range("D16").Value = 44711
$ws.Range("K16").Value = 8500
$ws.Range("L16").Value = 8500
$ws.Range("M16").Value = 8500
$ws.Range("O16").Value = 'Región Metropolitana'
$ws.Range("P16").Value = 236

# Row 17
$ws.Range("D17").Value = 44364
$ws.Range("J17").Value = 100

# Row 18
$ws.Range("D18").Value = 44362
$ws.Range("K18").Value = 6500
$ws.Range("L18").Value = 6500
$ws.Range("M18").Value = 6500
$ws.Range("N18").Value = '$/caja 36 atados'
$ws.Range("P18").Value = 181
$ws.Range("Q18").Value = 36

# Row 19
$ws.Range("D19").Value = 44371
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 6500
$ws.Range("L19").Value = 6500
$ws.Range("M19").Value = 6500
$ws.Range("P19").Value = 181

# Row 20
$ws.Range("D20").Value = 44354
$ws.Range("K20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 7000
$ws.Range("N20").Value = '$/caja 20 docenas'
$ws.Range("O20").Value = 'Región del Maule'
$ws.Range("P20").Value = 194
